$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("L15").Value = -14.285714285714
$ws.Range("N15").Value = -40
$ws.Range("C16").Value = 3
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 187.5
$ws.Range("I16").Value = 173
$ws.Range("J16").Value = 146
$ws.Range("K16").Value = 18.493150684931
$ws.Range("L16").Value = 76.530612244898
$ws.Range("M16").Value = -11.282051282051
$ws.Range("N16").Value = -55.297157622739
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 215
$ws.Range("J17").Value = 191
$ws.Range("K17").Value = 12.565445026178
$ws.Range("L17").Value = 9.693877551020
$ws.Range("M17").Value = 36.075949367088
$ws.Range("N17").Value = 10.824742268041
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 87
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = -6.451612903225
$ws.Range("L18").Value = 2.352941176470
$ws.Range("M18").Value = -66.409266409266
$ws.Range("N18").Value = -87.606837606837
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 30.952380952381
$ws.Range("I19").Value = 476
$ws.Range("J19").Value = 382
$ws.Range("K19").Value = 24.607329842931
$ws.Range("L19").Value = 31.491712707182
$ws.Range("M19").Value = 18.703241895261
$ws.Range("N19").Value = 30.410958904109
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -40
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 35
$ws.Range("H20").Value = -22.857142857142
$ws.Range("I20").Value = 253
$ws.Range("J20").Value = 255
$ws.Range("K20").Value = -0.784313725490
$ws.Range("L20").Value = 94.615384615384
$ws.Range("M20").Value = 42.134831460674
$ws.Range("N20").Value = -85.842193620593
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = -19.444444444444
$ws.Range("G21").Value = 109
$ws.Range("H21").Value = 18.348623853211
$ws.Range("I21").Value = 1220
$ws.Range("J21").Value = 1083
$ws.Range("K21").Value = 12.650046168051
$ws.Range("L21").Value = 37.542277339346
$ws.Range("M21").Value = 0.909842845326
$ws.Range("N21").Value = -64.821222606689
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = 8.333333333333
$ws.Range("L22").Value = -13.333333333333
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("I23").Value = 28
$ws.Range("J23").Value = 39
$ws.Range("K23").Value = -28.205128205128
$ws.Range("L23").Value = -17.647058823529
$ws.Range("M23").Value = -37.777777777777
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 181.818181818182
$ws.Range("F24").Value = 95
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = 15.853658536585
$ws.Range("I24").Value = 964
$ws.Range("J24").Value = 792
$ws.Range("K24").Value = 21.717171717171
$ws.Range("L24").Value = 26.013071895424
$ws.Range("M24").Value = -22.570281124498
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 40
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = 5.128205128205
$ws.Range("I25").Value = 383
$ws.Range("J25").Value = 370
$ws.Range("K25").Value = 3.513513513513
$ws.Range("L25").Value = 18.944099378882
$ws.Range("M25").Value = 5.801104972375
$ws.Range("L26").Value = 31.578947368421
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -20
$ws.Range("I27").Value = 46
$ws.Range("J27").Value = 37
$ws.Range("K27").Value = 24.324324324324
$ws.Range("L27").Value = 155.555555555556
$ws.Range("G28").Value = 2
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = 18.181818181818
$ws.Range("N28").Value = -31.578947368421
$ws.Range("G29").Value = 2
$ws.Range("J29").Value = 9
$ws.Range("K29").Value = 0
$ws.Range("N29").Value = -52.631578947368

$ws.Range("G15").Value = "'0"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").Value = "'***.*"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = 3
$ws.Range("I15").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = 0
$ws.Range("K15").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").Value = 1
$ws.Range("I15").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = -100
$ws.Range("K15").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Value = "'0"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null
$ws.Range("G22").Value = 1
$ws.Range("I15").Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4122) | Out-Null
$ws.Range("H22").Value = -100
$ws.Range("K15").Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value = 2
$ws.Range("I15").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("G26").Value = "'0"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("G26").PasteSpecial(-4122) | Out-Null
$ws.Range("H26").Value = "'***.*"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("H26").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = 1
$ws.Range("I15").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = -100
$ws.Range("K15").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").Value = 1
$ws.Range("I15").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = -100
$ws.Range("K15").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
